$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.111.37'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '3.028.88'
$ws.Range("E3").Value = '  +3.99%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '196.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '619.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.93%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.204'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.52%  '

$ws.Range("D10").Value = '3.026.65'
$ws.Range("E10").Value = '  +3.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.439'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.34%  '

$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.23'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.66%  '

$ws.Range("D14").Value = '3.588.87'
$ws.Range("E14").Value = '  +3.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.82'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.76%  '

$ws.Range("D16").Value = '76.085.70'
$ws.Range("E16").Value = '  +0.25%  '

$ws.Range("E17").Value = '  +2.21%  '

$ws.Range("D18").Value = '3.029.72'
$ws.Range("E18").Value = '  +3.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.44'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.40%  '

$ws.Range("E20").Value = '  +2.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.34%  '

$ws.Range("E22").Value = '  +5.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.12%  '

$ws.Range("D24").Value = '3.171.10'
$ws.Range("E24").Value = '  +3.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.71%  '

$ws.Range("E28").Value = '  +1.58%  '

$ws.Range("E29").Value = '  +2.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.03%  '

$ws.Range("E32").Value = '  +1.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '492.72'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.24%  '

$ws.Range("E34").Value = '  +5.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '20.53'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.74%  '

$ws.Range("E37").Value = '  +10.17%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '162.02'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.87%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '20.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.72%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '190.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.377'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.65%  '

$ws.Range("E42").Value = '  -5.25%  '

$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.770'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +18.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.14%  '

$ws.Range("E47").Value = '  +5.71%  '

$ws.Range("E48").Value = '  +0.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.16%  '

$ws.Range("E50").Value = '  +2.11%  '

$ws.Range("E51").Value = '  +0.73%  '
